$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-06 04:28:02"
$wsZhCn.Range("G5").Value = "2016-02-06 04:28:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-06 04:28:13"
$wsDeDe.Range("G5").Value = "2016-02-06 04:29:06"
